# Trade #2 closed at 2026-02-17 23:51:59 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate stats now that a 2nd trade has closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500.1   # Current Capital
$summary.Range("B4").Value = 0.1      # Total P&L $
$summary.Range("B5").Value = 1        # Total P&L %
$summary.Range("B6").Value = 2        # Total Trades
$summary.Range("B8").Value = 1        # Losing Trades
$summary.Range("B9").Value = 50       # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: update the MarketMaking strategy row (row 6).
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100.1     # Capital
$status.Range("D6").Value = 2         # Trades
$status.Range("E6").Value = 0.1       # P&L $
$status.Range("F6").Value = 0.1       # P&L %
$status.Range("G6").Value = 50        # Win Rate %

# ---------------------------------------------------------------------------
# Append the new closed trade (#2) to both "All Trades" and "MarketMaking".
# Date/time-looking text is force-written as text (number format "@") and
# then restyled back to Normal so it lands as plain text, not an auto
# -converted date/time serial, matching how the log values are stored.
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Add-TradeRow($ws) {
    $ws.Range("A3").Value = 2
    Set-TextValue $ws "B3" "2026-02-17"
    Set-TextValue $ws "C3" "23:51:53"
    $ws.Range("D3").Value = "MarketMaking"
    $ws.Range("E3").Value = "UP"
    $ws.Range("F3").Value = 0.92
    $ws.Range("G3").Value = 0.9
    $ws.Range("H3").Value = "CLOSED"
    $ws.Range("I3").Value = -2.1739
    $ws.Range("J3").Value = -0.02
    $ws.Range("K3").Value = 100.1
    $ws.Range("L3").Value = 0
    $ws.Range("M3").Value = 0
    $ws.Range("N3").Value = 0.6
    $ws.Range("O3").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P3").Value = "early_exit"
    $ws.Range("Q3").Value = 0.14
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
